# Update the "取得日時" (acquired datetime) column for all data rows
# on the "ランサーズ" sheet to reflect the new run timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-10 06:35:59"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
